$d = $word.ActiveDocument

# --- Locate the insertion point: the 'Tretaig hackspett' Heading1 paragraph in BILAGA 1 ---
$anchorRng = $d.Content
$found = $anchorRng.Find.Execute('Tretåig hackspett – ekologi samt krav på livsmiljön', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $found) { throw "Anchor heading not found" }
$anchorPar = $anchorRng.Paragraphs.First
$startIndex = $anchorPar.Index

# --- Build full plain-text block (13 paragraphs) and insert before anchor in one shot ---
$paraTexts = @(
    'Knärot – ekologi samt krav på livsmiljön',
    'Knärot är fridlyst enligt 8 och 15 §§ artskyddsförordningen och klassad som sårbar (VU) enligt rödlistan 2020. Knärot är beroende av hög och jämn luftfuktighet i gamla, ostörda skogsmiljöer och är känslig för snabba förändringar av ljus-/vindförhållanden eller uttorkning. På grund av ett alltför intensivt skogsbruk har den minskat med 40 (25-50) % under de senaste 60 åren och i framtiden bedöms minskningstakten uppgå till 30 (20-40) %. Till följd av att arten har en dokumenterat högre minskningstakt iförhållande till sin generationstid än vad som tidigare varit känt (data från Riksskogstaxeringen) höjdes den till hotkategori sårbar (VU) i rödlistan 2020 (Artdatabanken, 2021).',
    'Samuel Johnsons doktorsavhandling “Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“ (SLU, Uppsala 2014) visar att det krävs väl tilltagna skyddszoner för att knärotens växtplatser inte ska ta skada av skogsbruksåtgärder i intilliggande områden: “Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” Vidare “More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”',
    'Johnsons (2014) rekommendation på minst 50 meters breda skyddszoner runt knärotens växtplatser motsvarar en areal på 0,78 hektar, vilket ligger i linje med andra studier som gjorts på känsliga skogsarter: “In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”',
    'En nyligen publicerad vetenskaplig uppsats av Koelmeijer m.fl. (2022) inkluderar orkidén knärots skyddsbehov. I uppsatsen berörs problemet med uttorkning för växter, bl.a. för knärot, ett problem som blivit accentuerat på grund av den pågående klimatförändringen och torra somrar, t.ex. den exceptionellt torra sommaren 2018. I uppsatsen undersöks områden med tre olika avstånd från kalhyggeskant med avseende på skydd bl.a. för knärot. Det första området har avstånd upp till 20 m från hyggeskant (Strong edge effect), det andra 20 – 40 m från hyggeskant (Weak edge effect) och det tredje avser större avstånd från hyggeskant, där kanteffekten anses vara försumbar (Interior). Ett resultat var att man fann stor eller mycket stor uttorkningseffekt på känsliga och rödlistade skogsarter vid de kortare avstånden till hyggeskant, medan effekt av uttorkning inte konstaterades på större avstånd (Interior). För orkidén knärot fann man en rik förekomst (upp till 0,06 dm2/m2) på stort avstånd från hyggeskant (Interior), medan förekomsten var liten eller närmast försumbar i de områden som klassificerades som Weak edge effect respektive Strong edge effect. Arbetet påpekar att de allt oftare förekommande torra somrarna ger ytterligare skäl att utöka skyddsavståndet från hyggen till den fuktkrävande arten knärot (Koelmeijer m.fl., 2022).',
    'Även Skogsstyrelsens egen vägledning för hänsyn till knärot ligger i linje med ovanstående forskningsstudier. Av vägledningen framgår det att för med hög sannolikhet kunna bevara befintliga förekomster krävs relativt stora avsättningar av uppvuxen skog med slutet och relativt tätt kronskikt. Som riktlinje kan krävas ett avstånd på 50 meter in från brynet för att vidmakthålla ett fungerande mikroklimat. Detta innebär att fristående hänsynsytor för många arter (kärlväxter, lavar och mossor) kan behöva ha en area överstigande 0,8 hektar (cirkelyta med radien 50 meter = 0,78 hektar) för att bibehålla lokalklimatet. Även ganska små förändringar i form av förändrade ljus- och fuktighetsförhållanden, till exempel till följd av gallring, kan leda till att arten försvinner till följd av konkurrens med mera ljuskrävande och snabbväxande arter (Skogsstyrelsen, 2022).',
    'Referenser - knärot',
    'de Graaf M & Roberts M.R., 2009. Short-term response of the herbaceous layer within leave patches after harvest. Forest Ecology and Management 257, 1014-1025',
    'Halpern, C. B., Halaj, J., Evans, S. A., & Dovciak, M., 2012. Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. Ecological Applications, 22, 2049-2064 ',
    'Koelmeijer, I. A., Ehrlén, J., Jönsson, M., De Frenne, P., Berg, P., Andersson, J., Weibull, H. & Hylander, N. 2022. Interactive effects of drought and edge exposure on old-growth forest understory species. Landscape Ecology, 37, sid 1839-1853',
    'Rudolphi, J., Jönsson, M. T., & Gustafsson, L., 2014. Biological legacies buffer local species extinction after logging. Journal of Applied Ecology. 51, 53-62.',
    'Skogsstyrelsen, 2022. Vägledning för hänsyn till knärot. https://www.skogsstyrelsen.se/lag-och-tillsyn/artskydd/vagledningar-och-kunskapsstod-artskydd/vagledning-for-hansyn-till-knarot/',
    'SLU Artdatabanken, 2021. Artfaktablad. Naturvård – artfakta. SLU Artdatabanken, Uppsala '
)
$crChar = [char]13
$fullText = ($paraTexts -join $crChar) + $crChar
$anchorPar.Range.InsertBefore($fullText)

# --- Fix up paragraph styles for the newly inserted paragraphs ---
# Newly inserted paragraphs occupy indices [$startIndex .. $startIndex+12]; old heading now at $startIndex+13
$d.Paragraphs.Item($startIndex + 0).Style = 'Heading 1'
$d.Paragraphs.Item($startIndex + 1).Style = 'Normal'
$d.Paragraphs.Item($startIndex + 2).Style = 'Normal'
$d.Paragraphs.Item($startIndex + 3).Style = 'Normal'
$d.Paragraphs.Item($startIndex + 4).Style = 'Normal'
$d.Paragraphs.Item($startIndex + 5).Style = 'Normal'
$d.Paragraphs.Item($startIndex + 6).Style = 'Heading 2'
$d.Paragraphs.Item($startIndex + 7).Style = 'Normal'
$d.Paragraphs.Item($startIndex + 8).Style = 'Normal'
$d.Paragraphs.Item($startIndex + 9).Style = 'Normal'
$d.Paragraphs.Item($startIndex + 10).Style = 'Normal'
$d.Paragraphs.Item($startIndex + 11).Style = 'Normal'
$d.Paragraphs.Item($startIndex + 12).Style = 'Normal'

# --- Apply italics to specific runs, scoped to each paragraph's own range ---
$p3 = $d.Paragraphs.Item($startIndex + 2)
$r3_2 = $p3.Range.Duplicate
$ok3_2 = $r3_2.Find.Execute('“Retention Forestry as a Conservation Measure for Boreal Forest Ground Vegetation“', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok3_2) { throw "Italic run not found: paragraph 3 run 2" }
$r3_2.Font.Italic = 1
$r3_4 = $p3.Range.Duplicate
$ok3_4 = $r3_4.Find.Execute('“Study III shows that retention patches smaller than 0.5 ha do not lifeboat the sensitive forest herb G. repens, a species that depend on stable microclimatic conditions typical for intact forest stands.” ', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok3_4) { throw "Italic run not found: paragraph 3 run 4" }
$r3_4.Font.Italic = 1
$r3_6 = $p3.Range.Duplicate
$ok3_6 = $r3_6.Find.Execute('“More sensitive forest species are not lifeboated in retention patches ranging from 0.05 to 0.5 ha (Papers II & III).”', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok3_6) { throw "Italic run not found: paragraph 3 run 6" }
$r3_6.Font.Italic = 1

$p4 = $d.Paragraphs.Item($startIndex + 3)
$r4_2 = $p4.Range.Duplicate
$ok4_2 = $r4_2.Find.Execute('“In study III I also show that translocated specimens of G. repens survives well in mature forests at least 50 m from the nearest edge to an open area. Moreover, measures of temperature and humidity show that such distances from an open area is far enough to offer a microclimate that is more stable compared to what present in retention patches of around 0.1 ha. This means that the very centre of a circular patch with radius 50 m (equals a size of 0.78 ha) should offer conditions similar to interior forest and would perhaps be a suitable habitat for G. repens and similar species. Previous studies from both North America and Sweden have also concluded that patches between 0.5 and one ha are sufficient for preserving interior forest vegetation as well as sensitive lichens and bryophytes (de Graaf & Roberts 2009; Halpern et al. 2012; Rudolphi et al. 2014).”', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok4_2) { throw "Italic run not found: paragraph 4 run 2" }
$r4_2.Font.Italic = 1

$p8 = $d.Paragraphs.Item($startIndex + 7)
$r8_2 = $p8.Range.Duplicate
$ok8_2 = $r8_2.Find.Execute('Short-term response of the herbaceous layer within leave patches after harvest. ', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok8_2) { throw "Italic run not found: paragraph 8 run 2" }
$r8_2.Font.Italic = 1

$p9 = $d.Paragraphs.Item($startIndex + 8)
$r9_2 = $p9.Range.Duplicate
$ok9_2 = $r9_2.Find.Execute('Level and pattern of overstory retention interact to shape long-term responses of understories to timber harvest. ', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok9_2) { throw "Italic run not found: paragraph 9 run 2" }
$r9_2.Font.Italic = 1

$p10 = $d.Paragraphs.Item($startIndex + 9)
$r10_2 = $p10.Range.Duplicate
$ok10_2 = $r10_2.Find.Execute('Interactive effects of drought and edge exposure on old-growth forest understory species. ', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok10_2) { throw "Italic run not found: paragraph 10 run 2" }
$r10_2.Font.Italic = 1

$p11 = $d.Paragraphs.Item($startIndex + 10)
$r11_2 = $p11.Range.Duplicate
$ok11_2 = $r11_2.Find.Execute('Biological legacies buffer local species extinction after logging. ', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok11_2) { throw "Italic run not found: paragraph 11 run 2" }
$r11_2.Font.Italic = 1

$p12 = $d.Paragraphs.Item($startIndex + 11)
$r12_2 = $p12.Range.Duplicate
$ok12_2 = $r12_2.Find.Execute('Vägledning för hänsyn till knärot. ', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok12_2) { throw "Italic run not found: paragraph 12 run 2" }
$r12_2.Font.Italic = 1

$p13 = $d.Paragraphs.Item($startIndex + 12)
$r13_2 = $p13.Range.Duplicate
$ok13_2 = $r13_2.Find.Execute('Artfaktablad. Naturvård – artfakta. ', $true, $false, $false, $false, $false, $true, 1, $false, '', 0)
if (-not $ok13_2) { throw "Italic run not found: paragraph 13 run 2" }
$r13_2.Font.Italic = 1

# --- Update the header date from 2023-09-13 to 2023-09-15 ---
$sec = $d.Sections.First
$hdr = $sec.Headers.Item(2)
$dateOk = $hdr.Range.Find.Execute('2023-09-13', $true, $false, $false, $false, $false, $true, 1, $false, '2023-09-15', 2)
if (-not $dateOk) { throw "Date not found/replaced in header" }

Write-Output "DONE"
